{"js": "const replacements = [\n  [\"955\u00d72=\", \"230\u00d73=\"],\n  [\"977\u00d74=\", \"121\u00d73=\"],\n  [\"820\u00d77=\", \"452\u00d74=\"],\n  [\"553\u00d79=\", \"130\u00d75=\"],\n  [\"827\u00d74=\", \"349\u00d77=\"],\n  [\"482\u00d77=\", \"219\u00d73=\"],\n  [\"470\u00d74=\", \"507\u00d78=\"],\n  [\"375\u00d77=\", \"811\u00d77=\"],\n  [\"591\u00d79=\", \"745\u00d77=\"],\n  [\"430\u00d73=\", \"461\u00d73=\"],\n  [\"842\u00d72=\", \"142\u00d77=\"],\n  [\"975\u00d77=\", \"299\u00d73=\"],\n  [\"662\u00d73=\", \"492\u00d79=\"],\n  [\"579\u00d75=\", \"323\u00d73=\"],\n  [\"438\u00d74=\", \"501\u00d79=\"],\n  [\"914\u00d73=\", \"878\u00d72=\"],\n  [\"353\u00d77=\", \"326\u00d78=\"],\n  [\"456\u00d76=\", \"384\u00d74=\"],\n  [\"626\u00d76=\", \"985\u00d78=\"],\n  [\"250\u00d78=\", \"659\u00d77=\"],\n  [\"927\u00d74=\", \"459\u00d77=\"],\n  [\"266\u00d73=\", \"380\u00d75=\"],\n  [\"492\u00d73=\", \"846\u00d79=\"],\n  [\"877\u00d75=\", \"696\u00d76=\"],\n  [\"457\u00d77=\", \"423\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"955\u00d72=\", \"230\u00d73=\"),\n    @(\"977\u00d74=\", \"121\u00d73=\"),\n    @(\"820\u00d77=\", \"452\u00d74=\"),\n    @(\"553\u00d79=\", \"130\u00d75=\"),\n    @(\"827\u00d74=\", \"349\u00d77=\"),\n    @(\"482\u00d77=\", \"219\u00d73=\"),\n    @(\"470\u00d74=\", \"507\u00d78=\"),\n    @(\"375\u00d77=\", \"811\u00d77=\"),\n    @(\"591\u00d79=\", \"745\u00d77=\"),\n    @(\"430\u00d73=\", \"461\u00d73=\"),\n    @(\"842\u00d72=\", \"142\u00d77=\"),\n    @(\"975\u00d77=\", \"299\u00d73=\"),\n    @(\"662\u00d73=\", \"492\u00d79=\"),\n    @(\"579\u00d75=\", \"323\u00d73=\"),\n    @(\"438\u00d74=\", \"501\u00d79=\"),\n    @(\"914\u00d73=\", \"878\u00d72=\"),\n    @(\"353\u00d77=\", \"326\u00d78=\"),\n    @(\"456\u00d76=\", \"384\u00d74=\"),\n    @(\"626\u00d76=\", \"985\u00d78=\"),\n    @(\"250\u00d78=\", \"659\u00d77=\"),\n    @(\"927\u00d74=\", \"459\u00d77=\"),\n    @(\"266\u00d73=\", \"380\u00d75=\"),\n    @(\"492\u00d73=\", \"846\u00d79=\"),\n    @(\"877\u00d75=\", \"696\u00d76=\"),\n    @(\"457\u00d77=\", \"423\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 1\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
